$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F for rows 1-5 (header + first few data rows)
$ws.Range("F1").Value = "Material"
$ws.Range("F2").Value = "Carbonate"
$ws.Range("F3").Value = "Quartz"
$ws.Range("F4").Value = "Aragonite"
$ws.Range("F5").Value = "Dolomite"

# Row 20 gets its Material value next
$ws.Range("F20").Value = "Garnet"

# New row 21 (Zircon standard)
$ws.Range("A21").Value = "UWZ-1"
$ws.Range("D21").NumberFormat = $ws.Range("D1").NumberFormat

# New row 22 (Zircon standard)
$ws.Range("A22").Value = "KIM-5"
$ws.Range("C22").Value = 5.09
$ws.Range("D22").Value = "KIM\D*5"
$ws.Range("D22").NumberFormat = $ws.Range("D1").NumberFormat
$ws.Range("E22").Value = "Run"
$ws.Range("F22").Value = "Zircon"

$ws.Range("D21").Value = "UWZ\D*1"
$ws.Range("E21").Value = "Run"
$ws.Range("F21").Value = "Zircon"

# Fill remaining rows 6-19 of column F
$ws.Range("F6").Value = "Dolomite"
$ws.Range("F7").Value = "Ankerite"
$ws.Range("F8").Value = "Ankerite"
$ws.Range("F9").Value = "Ankerite"
$ws.Range("F10").Value = "Ankerite"
$ws.Range("F11").Value = "Ankerite"
$ws.Range("F12").Value = "Ankerite"
$ws.Range("F13").Value = "Ankerite"
$ws.Range("F14").Value = "Ankerite"
$ws.Range("F15").Value = "Ankerite"
$ws.Range("F16").Value = "Ankerite"
$ws.Range("F17").Value = "Ankerite"
$ws.Range("F18").Value = "Ankerite"
$ws.Range("F19").Value = "Ankerite"

$ws.Range("G24").Select()
